# Generate Report for Handback
# Updates the localization-status workbook to reflect that the zh-cn and
# de-de handoff files have now been handed back (in sync with en-US),
# filling in the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns and widening the columns that now
# hold longer text.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$zhHandbackFile = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$deHandbackFile = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$zhHandbackDate = "2016-09-01 10:43:18"
$deHandbackDate = "2016-09-01 10:43:25"
$urlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c377622dd272d2b852140f9c6dbc3431be9b4e2c/e2e/a.md"
$urlB = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c377622dd272d2b852140f9c6dbc3431be9b4e2c/e2e/b.md"

# ---------------------------------------------------------------------
# Overview sheet: "Ready for handoff" -> "Handed back: in sync with en-US"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# Widen the zh-cn / de-de status columns to fit the longer text.
$wsOverview.Columns.Item(5).ColumnWidth = 29.1
$wsOverview.Columns.Item(6).ColumnWidth = 29.1

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

# Latest Target File (I) + Latest Handback File (J)
$wsZh.Range("J2").Value = $zhHandbackFile
$wsZh.Range("J3").Value = $zhHandbackFile

# Latest Handback DateTime (K)
$wsZh.Range("K2").Value = $zhHandbackDate
$wsZh.Range("K3").Value = $zhHandbackDate

# Rebuild hyperlinks: A2, I2, A3, I3 (in that order) so the new "Latest
# Target File" links (column I) point at the same source doc as column A.
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $urlA, "", "", "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $urlA, "", "", "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $urlB, "", "", "b.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $urlA, "", "", "a.md")

# Widen the Status column (C) and Latest Handback File column (J).
$wsZh.Columns.Item(3).ColumnWidth = 29.1
$wsZh.Columns.Item(10).ColumnWidth = 39.1

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

# Latest Target File (I) + Latest Handback File (J)
$wsDe.Range("J2").Value = $deHandbackFile
$wsDe.Range("J3").Value = $deHandbackFile

# Latest Handback DateTime (K)
$wsDe.Range("K2").Value = $deHandbackDate
$wsDe.Range("K3").Value = $deHandbackDate

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $urlA, "", "", "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $urlA, "", "", "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $urlB, "", "", "b.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $urlA, "", "", "a.md")

$wsDe.Columns.Item(3).ColumnWidth = 29.1
$wsDe.Columns.Item(10).ColumnWidth = 39.1

Write-Host "Handback report generated."
